$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change FORMATO values from "pdf y png" to "pdf"
$ws.Range("J2").Value = "pdf"
$ws.Range("J3").Value = "pdf"

# Turn off Wrap Text for the cells that had it enabled
$wrapCells = @("E1","G1","H1","O1","P1","E2","G2","H2","O2","P2","E3","G3","H3","O3","P3")
foreach ($c in $wrapCells) {
    $ws.Range($c).WrapText = $false
}

# Rows auto-fit back to default height now that wrapping is off
$ws.Rows("1:3").AutoFit()

# Update selection
$ws.Range("L8").Select()
